$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 29
$ws.Cells.Item(3, 6).Value = 490
$ws.Cells.Item(5, 6).Value = 89
$ws.Cells.Item(6, 6).Value = 332
$ws.Cells.Item(7, 6).Value = 1313
$ws.Cells.Item(10, 6).Value = 1348
$ws.Cells.Item(11, 6).Value = 182
$ws.Cells.Item(13, 6).Value = 187
$ws.Cells.Item(14, 6).Value = 22
$ws.Cells.Item(15, 6).Value = 128
$ws.Cells.Item(16, 6).Value = 262
$ws.Cells.Item(17, 6).Value = 1700
$ws.Cells.Item(18, 6).Value = 629
$ws.Cells.Item(19, 6).Value = 273
$ws.Cells.Item(20, 6).Value = 286
$ws.Cells.Item(21, 6).Value = 3161
$ws.Cells.Item(22, 6).Value = 27
$ws.Cells.Item(23, 6).Value = 412
$ws.Cells.Item(24, 6).Value = 938
$ws.Cells.Item(25, 6).Value = 1220
$ws.Cells.Item(27, 6).Value = 2859
$ws.Cells.Item(28, 6).Value = 1657
$ws.Cells.Item(31, 6).Value = 675
$ws.Cells.Item(32, 6).Value = 874
$ws.Cells.Item(33, 6).Value = 5
$ws.Cells.Item(34, 6).Value = 1928
$ws.Cells.Item(35, 6).Value = 909
$ws.Cells.Item(36, 6).Value = 1932
$ws.Cells.Item(38, 6).Value = 361
$ws.Cells.Item(39, 6).Value = 83
$ws.Cells.Item(42, 6).Value = 911
$ws.Cells.Item(43, 6).Value = 816
$ws.Cells.Item(44, 6).Value = 1058
$ws.Cells.Item(45, 6).Value = 140
$ws.Cells.Item(46, 6).Value = 451
$ws.Cells.Item(47, 6).Value = 293
$ws.Cells.Item(48, 6).Value = 236
$ws.Cells.Item(49, 6).Value = 3382

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(8, 6).Value = 15
$ws.Cells.Item(13, 6).Value = 814

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 490
$ws.Cells.Item(4, 6).Value = 89
$ws.Cells.Item(7, 6).Value = 332
$ws.Cells.Item(8, 6).Value = 1313
$ws.Cells.Item(11, 6).Value = 1348
$ws.Cells.Item(12, 6).Value = 182
$ws.Cells.Item(14, 6).Value = 187
$ws.Cells.Item(15, 6).Value = 128
$ws.Cells.Item(16, 6).Value = 262
$ws.Cells.Item(17, 6).Value = 1700
$ws.Cells.Item(18, 6).Value = 629
$ws.Cells.Item(19, 6).Value = 273
$ws.Cells.Item(20, 6).Value = 286
$ws.Cells.Item(21, 6).Value = 3161
$ws.Cells.Item(22, 6).Value = 27
$ws.Cells.Item(23, 6).Value = 412
$ws.Cells.Item(24, 6).Value = 15
$ws.Cells.Item(25, 6).Value = 1220
$ws.Cells.Item(26, 6).Value = 2859
$ws.Cells.Item(27, 6).Value = 1657
$ws.Cells.Item(30, 6).Value = 814
$ws.Cells.Item(32, 6).Value = 874
$ws.Cells.Item(33, 6).Value = 1928
$ws.Cells.Item(35, 6).Value = 909
$ws.Cells.Item(36, 6).Value = 1932
$ws.Cells.Item(37, 6).Value = 361
$ws.Cells.Item(38, 6).Value = 83
$ws.Cells.Item(40, 6).Value = 911
$ws.Cells.Item(41, 6).Value = 816
$ws.Cells.Item(42, 6).Value = 1058
$ws.Cells.Item(43, 6).Value = 140
$ws.Cells.Item(44, 6).Value = 451
$ws.Cells.Item(45, 6).Value = 293
$ws.Cells.Item(47, 6).Value = 236
$ws.Cells.Item(48, 6).Value = 3382
